# Mise à jour de l'application
# Adds a new training-date column (AY) — 2025-09-22 (serial 45922) — right
# after the existing last column (AX, 2025-09-19), and records each
# player's attendance status ("P" = Présent, "B" = Blessure, "REP" = Repos)
# for that new date. Row 12 has no entry for the new date (left blank),
# matching the source edit. All dependent COUNTA/COUNTIF summary formulas
# in columns B:J recalc automatically once the new cells are written.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New date header for column AY (row 1), mirroring AX1's style/format.
# (Value is written before the format paste — writing the value AFTER a
# format-only PasteSpecial leaves the recalculation engine unaware that the
# COUNTA/COUNTIF summary formulas downstream need to re-evaluate.)
$ws.Range("AY1").Value = 45922
$ws.Range("AX1").Copy()
$ws.Range("AY1").PasteSpecial(-4122)   # xlPasteFormats

# Per-row attendance marker for the new date, mirroring each row's AX
# column style/format so the new cell matches existing ones exactly.
$attendance = [ordered]@{
  2  = "P"
  3  = "B"
  4  = "P"
  5  = "B"
  6  = "P"
  7  = "P"
  8  = "P"
  9  = "P"
  10 = "B"
  11 = "B"
  13 = "B"
  14 = "P"
  15 = "P"
  16 = "P"
  17 = "P"
  18 = "P"
  19 = "P"
  20 = "P"
  21 = "B"
  22 = "P"
  23 = "P"
  24 = "P"
  25 = "B"
  26 = "P"
  27 = "REP"
  28 = "P"
  29 = "P"
}

foreach ($row in $attendance.Keys) {
    $value = $attendance[$row]
    $srcCell = "AX" + $row
    $dstCell = "AY" + $row

    $ws.Range($dstCell).Value = $value
    $ws.Range($srcCell).Copy()
    $ws.Range($dstCell).PasteSpecial(-4122)   # xlPasteFormats
}

[void]($excel.CutCopyMode = $false)

# Mirror the author's post-edit selection/scroll nudge (one column over).
[void]$ws.Range("BA24").Select()
